$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row based on sheet data (data starts at row 2, header at row 1)
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
